$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[William%Zuccon%NULL%1,                Paolo%Comassi%NULL%1,                Luca%Adriani%NULL%1,                Giulio%Bergamaschini%NULL%1,                Elena%Bertin%NULL%1,                Raffaella%Borromeo%NULL%1,                Serena%Corti%NULL%1,                Federica%De Petri%NULL%1,                Francesco%Dolci%NULL%1,                Attilio%Galmozzi%NULL%1,                Alberto%Gigliotti%NULL%1,                Livio%Gualdoni%NULL%1,                Claudia%Guerra%NULL%1,                Anna%Khosthiova%NULL%1,                Giovanni%Leati%NULL%1,                Giuseppe%Lupi%NULL%1,                Paolo%Moscato%NULL%1,                Vittorio%Perotti%NULL%1,                Miriam%Piantelli%NULL%1,                Alain%Ruini%NULL%1,                Silvia%Sportelli%NULL%1,                Micaela%Susca%NULL%1,                Carmine%Troiano%NULL%1,                Giampaolo%Benelli%NULL%1,                Elisabetta%Buscarini%NULL%1,                Ciro%Canetta%NULL%1,                Guido%Merli%NULL%1,                Alessandro%Scartabellati%NULL%1]'
$ws.Range("E3").Value2 = '[Imen Ben%Saida%NULL%1,                Emna%Ennouri%NULL%1,                Rayane%Nachi%NULL%1,                Khaoula%Meddeb%NULL%1,                Jihene%Mahmoud%NULL%1,                Nesrine%Thabet%NULL%1,                Salma%Jerbi%NULL%1,                Mohamed%Boussarsar%NULL%1]'
$ws.Range("E4").Value2 = '[Ilias I.%Siempos%NULL%1,                Eleni%Xourgia%NULL%1,                Theodora K.%Ntaidou%NULL%1,                Dimitris%Zervakis%NULL%1,                Eleni E.%Magira%NULL%1,                Anastasia%Kotanidou%NULL%1,                Christina%Routsi%NULL%1,                Spyros G.%Zakynthinos%NULL%1]'
$ws.Range("E5").Value2 = '[Christian%Karagiannidis%NULL%1,                Carina%Mostert%NULL%1,                Corinna%Hentschker%NULL%1,                Thomas%Voshaar%NULL%1,                Jürgen%Malzahn%NULL%1,                Gerhard%Schillinger%NULL%1,                Jürgen%Klauber%NULL%1,                Uwe%Janssens%NULL%1,                Gernot%Marx%NULL%1,                Steffen%Weber-Carstens%NULL%1,                Stefan%Kluge%NULL%0,                Michael%Pfeifer%NULL%1,                Linus%Grabenhenrich%NULL%1,                Tobias%Welte%NULL%1,                Reinhard%Busse%NULL%1]'
$ws.Range("E6").Value2 = '[Matthieu%Schmidt%NULL%1,                David%Hajage%NULL%2,                David%Hajage%NULL%0,                Alexandre%Demoule%NULL%1,                Tài%Pham%NULL%1,                Alain%Combes%NULL%1,                Martin%Dres%NULL%1,                Said%Lebbah%NULL%1,                Antoine%Kimmoun%NULL%1,                Alain%Mercat%NULL%1,                Gaëtan%Beduneau%NULL%1,                Jessica%Palmyre%NULL%1,                Margot%Prevost%NULL%1,                Pierre%Asfar%NULL%1,                François%Beloncle%NULL%1,                Julien%Demiselle%NULL%1,                Arthur%Pavot%NULL%1,                Xavier%Monnet%NULL%1,                Christian%Richard%NULL%1,                Julien%Mayaux%NULL%1,                Alexandra%Beurton%NULL%1,                Richard%Descamps%NULL%1,                Aurélie%Joret%NULL%1,                Damien%Du Cheyron%NULL%1,                Frédéric%Pene%NULL%1,                Jean-Daniel%Chiche%NULL%1,                Mathieu%Jozwiak%NULL%1,                Paul%Jaubert%NULL%1,                Guillaume%Voiriot%NULL%1,                Muriel%Fartoukh%NULL%1,                Marion%Teulier%NULL%1,                Clarisse%Blayau%NULL%1,                Erwen%L’Her%NULL%1,                Cécile%Aubron%NULL%1,                Laetitia%Bodenes%NULL%1,                Nicolas%Ferriere%NULL%1,                Johann%Auchabie%NULL%1,                Anthony%Le Meur%NULL%1,                Sylvain%Pignal%NULL%1,                Thierry%Mazzoni%NULL%1,                Jean-Pierre%Quenot%NULL%1,                Pascal%Andreu%NULL%1,                Jean-Baptiste%Roudau%NULL%1,                Marie%Labruyère%NULL%1,                Saad%Nseir%NULL%0,                Sébastien%Preau%NULL%1,                Julien%Poissy%NULL%0,                Daniel%Mathieu%NULL%0,                Sarah%Benhamida%NULL%1,                Rémi%Paulet%NULL%1,                Nicolas%Roucaud%NULL%1,                Martial%Thyrault%NULL%1,                Florence%Daviet%NULL%1,                Sami%Hraiech%NULL%1,                Gabriel%Parzy%NULL%1,                Aude%Sylvestre%NULL%1,                Sébastien%Jochmans%NULL%1,                Anne-Laure%Bouilland%NULL%1,                Mehran%Monchi%NULL%1,                Marc%Danguy des Déserts%NULL%1,                Quentin%Mathais%NULL%1,                Gwendoline%Rager%NULL%1,                Pierre%Pasquier%NULL%1,                Reignier%Jean%NULL%1,                Seguin%Amélie%NULL%1,                Garret%Charlotte%NULL%1,                Canet%Emmanuel%NULL%1,                Jean%Dellamonica%NULL%1,                Clément%Saccheri%NULL%1,                Romain%Lombardi%NULL%1,                Yanis%Kouchit%NULL%1,                Sophie%Jacquier%NULL%1,                Armelle%Mathonnet%NULL%1,                Mai-Ahn%Nay%NULL%1,                Isabelle%Runge%NULL%1,                Frédéric%Martino%NULL%1,                Laure%Flurin%NULL%1,                Amélie%Rolle%NULL%1,                Michel%Carles%NULL%1,                Rémi%Coudroy%NULL%1,                Arnaud W.%Thille%NULL%1,                Jean-Pierre%Frat%NULL%1,                Maeva%Rodriguez%NULL%1,                Pascal%Beuret%NULL%1,                Audrey%Tientcheu%NULL%1,                Arthur%Vincent%NULL%1,                Florian%Michelin%NULL%1,                Marie Anne%Melone%NULL%1,                Maxime%Gauzi%NULL%1,                Arnaud%Guilbert%NULL%1,                Geoffrey%Kouadri%NULL%1,                Valérie%Gissot%NULL%1,                Stéphan%Ehrmann%NULL%1,                Charlotte%Salmon-Gandonnière%NULL%1,                Djlali%Elaroussi%NULL%1,                Agathe%Delbove%NULL%1,                Yannick%Fedun%NULL%1,                Julien%Huntzinger%NULL%1,                Eddy%Lebas%NULL%1,                Grâce%Kisoka%NULL%1,                Céline%Grégoire%NULL%1,                Stella%Marchetta%NULL%1,                Bernard%Lambermont%NULL%1,                Laurent%Argaud%NULL%1,                Thomas%Baudry%NULL%1,                Pierre-Jean%Bertrand%NULL%1,                Auguste%Dargent%NULL%1,                Christophe%Guitton%NULL%1,                Nicolas%Chudeau%NULL%1,                Mickaël%Landais%NULL%1,                Cédric%Darreau%NULL%1,                Alexis%Ferre%NULL%1,                Antoine%Gros%NULL%1,                Guillaume%Lacave%NULL%1,                Fabrice%Bruneel%NULL%1,                Mathilde%Neuville%NULL%1,                Jérôme%Devaquet%NULL%1,                Guillaume%Tachon%NULL%1,                Richard%Gallot%NULL%1,                Riad%Chelha%NULL%1,                Arnaud%Galbois%NULL%1,                Anne%Jallot%NULL%1,                Ludivine Chalumeau%Lemoine%NULL%1,                Khaldoun%Kuteifan%NULL%1,                Valentin%Pointurier%NULL%1,                Louise-Marie%Jandeaux%NULL%1,                Joy%Mootien%NULL%1,                Charles%Damoisel%NULL%1,                Benjamin%Sztrymf%NULL%1,                Juliette%Chommeloux%NULL%1,                Charles Edouard%Luyt%NULL%1,                Frédérique%Schortgen%NULL%1,                Leon%Rusel%NULL%1,                Camille%Jung%NULL%1,                Florent%Gobert%NULL%1,                Damien%Vimpere%NULL%1,                Lionel%Lamhaut%NULL%1,                Bertrand%Sauneuf%NULL%1,                Liliane%Charrier%NULL%1,                Julien%Calus%NULL%1,                Isabelle%Desmeules%NULL%1,                Benoît%Painvin%NULL%1,                Jean-Marc%Tadie%NULL%1,                Vincent%Castelain%NULL%1,                Baptiste%Michard%NULL%1,                Jean-Etienne%Herbrecht%NULL%1,                Mathieu%Baldacini%NULL%1,                Nicolas%Weiss%NULL%1,                Sophie%Demeret%NULL%1,                Clémence%Marois%NULL%1,                Benjamin%Rohaut%NULL%1,                Pierre-Henri%Moury%NULL%1,                Anne-Charlotte%Savida%NULL%1,                Emmanuel%Couadau%NULL%1,                Mathieu%Série%NULL%1,                Nica%Alexandru%NULL%1,                Cédric%Bruel%NULL%1,                Candice%Fontaine%NULL%1,                Sonia%Garrigou%NULL%1,                Juliette Courtiade%Mahler%NULL%1,                Maxime%Leclerc%NULL%1,                Michel%Ramakers%NULL%1,                Pierre%Garçon%NULL%1,                Nicole%Massou%NULL%1,                Ly Phacs%Van Vong%NULL%1,                Juliane%Sen%NULL%1,                Nolwenn%Lucas%NULL%1,                Franck%Chemouni%NULL%1,                Annabelle%Stoclin%NULL%1,                Alexandre%Avenel%NULL%1,                Henri%Faure%NULL%1,                Angélie%Gentilhomme%NULL%1,                Sylvie%Ricome%NULL%1,                Paul%Abraham%NULL%1,                Céline%Monard%NULL%1,                Julien%Textoris%NULL%1,                Thomas%Rimmele%NULL%1,                Florent%Montini%NULL%1,                Gabriel%Lejour%NULL%1,                Thierry%Lazard%NULL%1,                Isabelle%Etienney%NULL%1,                Younes%Kerroumi%NULL%1,                Dupuis%Claire%NULL%1,                Bereiziat Marine%Coupez%NULL%1,                Thouy%François%NULL%1,                Clémet%Hoffmann%NULL%1,                Nicolas%Donat%NULL%1,                Violaine%Muller%NULL%1,                Thibault%Martinez%NULL%1,                Audrey%Jacquot%NULL%1,                Matthieu%Mattei%NULL%1,                Bruno%Levy%NULL%1,                Ramin%Ravan%NULL%1,                Loïc%Dopeux%NULL%1,                Jean-Mathias%Liteaudon%NULL%1,                Delphine%Roux%NULL%1,                Brice%Rey%NULL%1,                Radu%Anghel%NULL%1,                Deborah%Schenesse%NULL%1,                Vincent%Gevrey%NULL%1,                Jermy%Castanera%NULL%1,                Philippe%Petua%NULL%1,                Benjamin%Madeux%NULL%1,                Otto%Hartman%NULL%1,                Michael%Piagnerelli%NULL%1,                Anne%Joosten%NULL%1,                Cinderella%Noel%NULL%1,                Patrick%Biston%NULL%1,                Thibaut%Noel%NULL%1,                Gurvan L. E.%Bouar%NULL%1,                Messabi%Boukhanza%NULL%1,                Elsa%Demarest%NULL%1,                Marie-France%Bajolet%NULL%1,                Nathanaël%Charrier%NULL%1,                Audrey%Quenet%NULL%1,                Cécile%Zylberfajn%NULL%1,                Nicolas%Dufour%NULL%1,                Buno%Mégarbane%NULL%1,                Sqébastian%Voicu%NULL%1,                Nicolas%Deye%NULL%1,                Isabelle%Malissin%NULL%1,                François%Legay%NULL%1,                Matthieu%Debarre%NULL%1,                Nicolas%Barbarot%NULL%1,                Pierre%Fillatre%NULL%1,                Bertrand%Delord%NULL%1,                Thomas%Laterrade%NULL%1,                Tahar%Saghi%NULL%1,                Wilfried%Pujol%NULL%1,                Pierre Julien%Cungi%NULL%1,                Pierre%Esnault%NULL%1,                Mickael%Cardinale%NULL%1,                Vivien%Hong Tuan Ha%NULL%1,                Grégory%Fleury%NULL%1,                Marie-Ange%Brou%NULL%1,                Daniel%Zafimahazo%NULL%1,                David%Tran-Van%NULL%1,                Patrick%Avargues%NULL%1,                Lisa%Carenco%NULL%1,                Nicolas%Robin%NULL%1,                Alexandre%Ouali%NULL%1,                Lucie%Houdou%NULL%1,                Christophe%Le Terrier%NULL%1,                Noémie%Suh%NULL%1,                Steve%Primmaz%NULL%1,                Jérome%Pugin%NULL%1,                Emmanuel%Weiss%NULL%1,                Tobias%Gauss%NULL%1,                Jean-Denis%Moyer%NULL%1,                Catherine%Paugam-Burtz%NULL%1,                Béatrice%La Combe%NULL%1,                Rolland%Smonig%NULL%1,                Jade%Violleau%NULL%1,                Pauline%Cailliez%NULL%1,                Jonathan%Chelly%NULL%1,                Antoine%Marchalot%NULL%1,                Cécile%Saladin%NULL%1,                Christelle%Bigot%NULL%1,                Pierre-Marie%Fayolle%NULL%1,                Jules%Fatséas%NULL%1,                Amr%Ibrahim%NULL%1,                Dabor%Resiere%NULL%1,                Rabih%Hage%NULL%1,                Clémentine%Cholet%NULL%1,                Marie%Cantier%NULL%1,                Pierre%Trouiller%NULL%1,                Philippe%Montravers%NULL%1,                Brice%Lortat-Jacob%NULL%1,                Sebastien%Tanaka%NULL%1,                Alexy%Tran-Dinh%NULL%1,                Jacques%Duranteau%NULL%1,                Anatole%Harrois%NULL%1,                Guillaume%Dubreuil%NULL%1,                Marie%Werner%NULL%1,                Anne%Godier%NULL%1,                Sophie%Hamada%NULL%1,                Diane%Zlotnik%NULL%1,                Hélène%Nougue%NULL%1,                Armand%Mekontso-Dessap%NULL%1,                Guillaume%Carteaux%NULL%1,                Keyvan%Razazi%NULL%1,                Nicolas%De Prost%NULL%1,                Nicolas%Mongardon%NULL%1,                Olivier%Langeron%NULL%1,                Eric%Levesque%NULL%1,                Arié%Attias%NULL%1,                Charles%de Roquetaillade%NULL%1,                Benjamin G.%Chousterman%NULL%1,                Alexandre%Mebazaa%NULL%1,                Etienne%Gayat%NULL%1,                Marc%Garnier%NULL%1,                Emmanuel%Pardo%NULL%1,                Lea%Satre-Buisson%NULL%1,                Christophe%Gutton%NULL%1,                Elise%Yvin%NULL%1,                Clémence%Marcault%NULL%1,                Elie%Azoulay%NULL%1,                Michael%Darmon%NULL%1,                Nicolas%Bonnet%NULL%1,                Nathan%Ebstein%NULL%1,                Stéphane%Gaudry%NULL%1,                Yves%Cohen%NULL%1,                 Hafid%Ait-Oufella%NULL%1,                Geoffroy%Hariri%NULL%1,                Tomas%Urbina%NULL%1,                Sandie%Mazerand%NULL%1,                Nicholas%Heming%NULL%1,                Francesca%Santi%NULL%1,                Pierre%Moine%NULL%1,                Djillali%Annane%NULL%1,                Adrien%Bouglé%NULL%1,                Edris%Omar%NULL%1,                Aymeric%Lancelot%NULL%1,                Emmanuelle%Begot%NULL%1,                Gaétan%Plantefeve%NULL%1,                Damien%Contou%NULL%1,                Hervé%Mentec%NULL%1,                Olivier%Pajot%NULL%1,                Stanislas%Faguer%NULL%1,                Olivier%Cointault%NULL%1,                Laurence%Lavayssiere%NULL%1,                Marie-Béatrice%Nogier%NULL%1,                Matthieu%Jamme%NULL%1,                Claire%Pichereau%NULL%1,                Jan%Hayon%NULL%1,                Hervé%Outin%NULL%1,                François%Dépret%NULL%1,                Maxime%Coutrot%NULL%1,                Maité%Chaussard%NULL%1,                Lucie%Guillemet%NULL%1,                Pierre%Goffin%NULL%1,                Romain%Thouny%NULL%1,                Julien%Guntz%NULL%1,                Laurent%Jadot%NULL%1,                Romain%Persichini%NULL%1,                Vanessa%Jean-Michel%NULL%1,                Hugues%Georges%NULL%1,                Thomas%Caulier%NULL%1,                Gaël%Pradel%NULL%1,                Marie-Hélène%Hausermann%NULL%1,                ThiMy Hue%Nguyen-Valat%NULL%1,                Michel%Boudinaud%NULL%1,                Emmanuel%Vivier%NULL%1,                Sylvène%Rosseli%NULL%1,                Gaël%Bourdin%NULL%1,                Christian%Pommier%NULL%1,                Marc%Vinclair%NULL%1,                Simon%Poignant%NULL%1,                Sandrine%Mons%NULL%1,                Wulfran%Bougouin%NULL%1,                Franklin%Bruna%NULL%1,                Quentin%Maestraggi%NULL%1,                Christian%Roth%NULL%1,                Laurent%Bitker%NULL%1,                François%Dhelft%NULL%1,                Justine%Bonnet-Chateau%NULL%1,                Mathilde%Filippelli%NULL%1,                Tristan%Morichau-Beauchant%NULL%1,                Stéphane%Thierry%NULL%1,                Charlotte%Le Roy%NULL%1,                Mélanie%Saint Jouan%NULL%1,                Bruno%Goncalves%NULL%1,                Aurélien%Mazeraud%NULL%1,                Matthieu%Daniel%NULL%1,                Tarek%Sharshar%NULL%1,                Cyril%Cadoz%NULL%1,                Rostane%Gaci%NULL%1,                Sébastien%Gette%NULL%1,                Guillaune%Louis%NULL%1,                Sophe-Caroline%Sacleux%NULL%1,                Marie-Amélie%Ordan%NULL%1,                Aurélie%Cravoisy%NULL%1,                Marie%Conrad%NULL%1,                Guilhem%Courte%NULL%1,                Sébastien%Gibot%NULL%1,                Younès%Benzidi%NULL%1,                Claudia%Casella%NULL%1,                Laurent%Serpin%NULL%1,                Jean-Lou%Setti%NULL%1,                Marie-Catherine%Besse%NULL%1,                Anna%Bourreau%NULL%1,                Jérôme%Pillot%NULL%1,                Caroline%Rivera%NULL%1,                Camille%Vinclair%NULL%1,                Marie-Aline%Robaux%NULL%1,                Chloé%Achino%NULL%1,                Marie-Charlotte%Delignette%NULL%1,                Tessa%Mazard%NULL%1,                Frédéric%Aubrun%NULL%1,                Bruno%Bouchet%NULL%1,                Aurélien%Frérou%NULL%1,                Laura%Muller%NULL%1,                Charlotte%Quentin%NULL%1,                Samuel%Degoul%NULL%1,                Xavier%Stihle%NULL%1,                Claude%Sumian%NULL%1,                Nicoletta%Bergero%NULL%1,                Bernard%Lanaspre%NULL%1,                Hervé%Quintard%NULL%1,                Eve Marie%Maiziere%NULL%1,                Pierre-Yves%Egreteau%NULL%1,                Guillaume%Leloup%NULL%1,                Florin%Berteau%NULL%1,                Marjolaine%Cottrel%NULL%1,                Marie%Bouteloup%NULL%1,                Matthieu%Jeannot%NULL%1,                Quentin%Blanc%NULL%1,                Julien%Saison%NULL%1,                Isabelle%Geneau%NULL%1,                Romaric%Grenot%NULL%1,                Abdel%Ouchike%NULL%1,                Pascal%Hazera%NULL%1,                Anne-Lyse%Masse%NULL%1,                Suela%Demiri%NULL%1,                Corinne%Vezinet%NULL%1,                Elodie%Baron%NULL%1,                Déborah%Benchetrit%NULL%1,                Antoine%Monsel%NULL%1,                Grégoire%Trebbia%NULL%1,                Emmanuelle%Schaack%NULL%1,                Raphaël%Lepecq%NULL%1,                Mathieu%Bobet%NULL%1,                Christophe%Vinsonneau%NULL%1,                Thibault%Dekeyser%NULL%1,                Quentin%Delforge%NULL%1,                Imen%Rahmani%NULL%1,                Bérengère%Vivet%NULL%1,                Jonathan%Paillot%NULL%1,                Lucie%Hierle%NULL%1,                Claire%Chaignat%NULL%1,                Sarah%Valette%NULL%1,                Benoït%Her%NULL%1,                Jennifier%Brunet%NULL%1,                Mathieu%Page%NULL%1,                Fabienne%Boiste%NULL%1,                Anthony%Collin%NULL%1,                Florent%Bavozet%NULL%2,                Aude%Garin%NULL%1,                Mohamed%Dlala%NULL%1,                Kais%Mhamdi%NULL%1,                Bassem%Beilouny%NULL%1,                Alexandra%Lavalard%NULL%1,                Severine%Perez%NULL%1,                Benoit%Veber%NULL%1,                Pierre-Gildas%Guitard%NULL%1,                Philippe%Gouin%NULL%1,                Anna%Lamacz%NULL%1,                Fabienne%Plouvier%NULL%1,                Bertrand P.%Delaborde%NULL%1,                Aïssa%Kherchache%NULL%1,                Amina%Chaalal%NULL%1,                Jean-Damien%Ricard%NULL%1,                Marc%Amouretti%NULL%1,                Santiago%Freita-Ramos%NULL%1,                Damien%Roux%NULL%1,                Jean-Michel%Constantin%NULL%1,                Mona%Assefi%NULL%1,                Marine%Lecore%NULL%1,                Agathe%Selves%NULL%1,                Florian%Prevost%NULL%1,                Christian%Lamer%NULL%1,                Ruiying%Shi%NULL%1,                Lyes%Knani%NULL%1,                Sébastien%Pili-Floury%NULL%1,                Lucie%Vettoretti%NULL%1,                Michael%Levy%NULL%1,                Lucile%Marsac%NULL%1,                Stéphane%Dauger%NULL%1,                Sophie%Guilmin-Crépon%NULL%1,                Jean-Baptiste%Putegnat%NULL%1,                Frédérique%Bayle%NULL%1,                Maya%Perrou%NULL%1,                Ghyslaine%Thao%NULL%1,                Guillaume%Géri%NULL%1,                Cyril%Charron%NULL%1,                Xavier%Repessé%NULL%1,                Antoine%Vieillard-Baron%NULL%1,                Mathieu%Guilbart%NULL%1,                Pierre-Alexandre%Roger%NULL%1,                Sébastien%Hinard%NULL%1,                Pierre-Yves%Macq%NULL%1,                Kevin%Chaulier%NULL%1,                Sylvie%Goutte%NULL%1,                Patrick%Chillet%NULL%1,                Anaïs%Pitta%NULL%1,                Barbara%Darjent%NULL%1,                Amandine%Bruneau%NULL%1,                Sigismond%Lasocki%NULL%1,                Maxime%Leger%NULL%1,                Soizic%Gergaud%NULL%1,                Pierre%Lemarie%NULL%1,                Nicolas%Terzi%NULL%1,                Carole%Schwebel%NULL%1,                Anaïs%Dartevel%NULL%1,                Louis-Marie%Galerneau%NULL%1,                Jean-Luc%Diehl%NULL%1,                Caroline%Hauw-Berlemont%NULL%0,                Nicolas%Péron%NULL%1,                Emmanuel%Guérot%NULL%1,                Abolfazl Mohebbi%Amoli%NULL%1,                Michel%Benhamou%NULL%1,                Jean-Pierre%Deyme%NULL%1,                Olivier%Andremont%NULL%1,                Diane%Lena%NULL%1,                Julien%Cady%NULL%1,                Arnaud%Causeret%NULL%1,                Arnaud%De La Chapelle%NULL%1,                Christophe%Cracco%NULL%1,                Stéphane%Rouleau%NULL%1,                David%Schnell%NULL%1,                Cécile%Lory%NULL%1,                Thibault%Chapelle%NULL%1,                Vincent%Bruckert%NULL%1,                Julie%Garcia%NULL%1,                Abdlazize%Sahraoui%NULL%1,                Nathalie%Abbosh%NULL%1,                Caroline%Bornstain%NULL%1,                Pierre%Pernet%NULL%1,                Florent%Poirson%NULL%1,                Ahmed%Pasem%NULL%1,                Philippe%Karoubi%NULL%1,                Virginie%Poupinel%NULL%1,                Caroline%Gauthier%NULL%1,                François%Bouniol%NULL%1,                Philippe%Feuchere%NULL%1,                Florent%Bavozet%NULL%0,                Anne%Heron%NULL%1,                Serge%Carreira%NULL%1,                Malo%Emery%NULL%1,                Anne%Le Floch%NULL%1,                Luana%Giovannangeli%NULL%1,                Nicolas%Herzog%NULL%1,                Christophe%Giacardi%NULL%1,                Thibaut%Baudic%NULL%1,                Chloé%Thill%NULL%1,                Florence%Tubach%NULL%1,                Olivier%Lesieur%NULL%1,                Julie%Noublanche%NULL%1]'
$ws.Range("E7").Value2 = '[Alfonso C.%Hernandez-Romieu%NULL%1,                Max W.%Adelman%NULL%2,                Maxwell A.%Hockstein%NULL%1,                Chad J.%Robichaux%NULL%1,                Johnathan A.%Edwards%NULL%1,                Jane C.%Fazio%NULL%1,                James M.%Blum%NULL%3,                Craig S.%Jabaley%NULL%2,                Mark%Caridi-Scheible%NULL%2,                Greg S.%Martin%NULL%2,                David J.%Murphy%NULL%3,                Sara C.%Auld%NULL%1]'
$ws.Range("E8").Value2 = '[Albano%G.%coreGivesNoEmail%0,              Antonelli%M.%coreGivesNoEmail%0,              Bellani%G.%coreGivesNoEmail%0,              Bonanomi%E.%coreGivesNoEmail%0,              Cabrini%L.%coreGivesNoEmail%0,              Carlesso%E.%coreGivesNoEmail%0,              Castelli%G.%coreGivesNoEmail%0,              Cattaneo%S.%coreGivesNoEmail%0,              Cecconi%M.%coreGivesNoEmail%0,              Cereda%D.%coreGivesNoEmail%0,              Colombo%S.%coreGivesNoEmail%0,              Coluccello%A.%coreGivesNoEmail%0,              Crescini%G.%coreGivesNoEmail%0,              Forastieri%Molinari A.%coreGivesNoEmail%0,              Foti%G.%coreGivesNoEmail%0,              Fumagalli%R.%coreGivesNoEmail%0,              Grasselli%G.%coreGivesNoEmail%0,              Greco%M.%coreGivesNoEmail%0,              Iotti%G. A.%coreGivesNoEmail%0,              Langer%T.%coreGivesNoEmail%0,              Latronico%N.%coreGivesNoEmail%0,              Lorini%F. L.%coreGivesNoEmail%0,              Mojoli%F.%coreGivesNoEmail%0,              Natalini%G.%coreGivesNoEmail%0,              Pesenti%A.%coreGivesNoEmail%0,              Pessina%C. M.%coreGivesNoEmail%0,              Ranieri%V. M.%coreGivesNoEmail%0,              Rech%R.%coreGivesNoEmail%0,              Rosano%A.%coreGivesNoEmail%0,              Scudeller%L.%coreGivesNoEmail%0,              Storti%E.%coreGivesNoEmail%0,              Thompson%B. T.%coreGivesNoEmail%0,              Tirani%M.%coreGivesNoEmail%0,              Villani%P. G.%coreGivesNoEmail%0,              Zanella%A.%coreGivesNoEmail%0]'
$ws.Range("E9").Value2 = '[Yong Hoon%Lee%NULL%1,                Keum-Ju%Choi%NULL%1,                Sun Ha%Choi%NULL%1,                Shin Yup%Lee%NULL%1,                Kyung Chan%Kim%NULL%1,                Eun Jin%Kim%NULL%1,                Jaehee%Lee%NULL%2,                Jaehee%Lee%NULL%0]'
$ws.Range("E10").Value2 = '[Atul%Matta%NULL%1,                Siddique%Chaudhary%NULL%1,                Kevin%Bryan Lo%NULL%2,                Kevin%Bryan Lo%NULL%0,                Robert%DeJoy%NULL%1,                Fahad%Gul%NULL%1,                Ricardo%Torres%NULL%1,                Neal%Chaisson%NULL%2,                Neal%Chaisson%NULL%0,                Gabriel%Patarroyo-Aponte%NULL%2,                Gabriel%Patarroyo-Aponte%NULL%0]'
$ws.Range("E11").Value2 = '[Ricard%Mellado-Artigas%rmartigas@gmail.com%1,                Bruno L.%Ferreyro%NULL%2,                Bruno L.%Ferreyro%NULL%0,                Federico%Angriman%NULL%1,                María%Hernández-Sanz%NULL%2,                Egoitz%Arruti%NULL%1,                Antoni%Torres%NULL%2,                Jesús%Villar%NULL%1,                Laurent%Brochard%NULL%1,                Carlos%Ferrando%NULL%2,                Ricard%Mellado-Artigas%NULL%1,                María%Hernández-Sanz%NULL%0,                Carlos%Ferrando%NULL%0,                Marina%Vendrell%NULL%1,                Gerard%Sánchez-Etayo%NULL%1,                Amalia%Alcón%NULL%1,                Isabel%Belda%NULL%1,                Mercé%Agustí%NULL%1,                Albert%Carramiñana%NULL%1,                Isabel%Gracia%NULL%1,                Miriam%Panzeri%NULL%1,                Irene%León%NULL%1,                Jaume%Balust%NULL%1,                Ricard%Navarro%NULL%2,                María José%Arguís%NULL%1,                María José%Carretero%NULL%1,                Cristina%Ibáñez%NULL%1,                Juan%Perdomo%NULL%1,                Antonio%López%NULL%1,                Manuel%López-Baamonde%NULL%1,                Tomás%Cuñat%NULL%1,                Marta%Ubré%NULL%1,                Antonio%Ojeda%NULL%1,                Andrea%Calvo%NULL%1,                Eva%Rivas%NULL%1,                Paola%Hurtado%NULL%1,                Roger%Pujol%NULL%1,                Nuria%Martín%NULL%1,                Javier%Tercero%NULL%1,                Pepe%Sanahuja%NULL%1,                Marta%Magaldi%NULL%1,                Miquel%Coca%NULL%1,                Elena%del Rio%NULL%1,                Julia%Martínez-Ocon%NULL%1,                Paula%Masgoret%NULL%1,                Monserrat%Tio%NULL%1,                Angel%Caballero%NULL%1,                Raquel%Risco%NULL%1,                Raquel%Bergé%NULL%1,                Lidia%Gómez%NULL%1,                Nicolás%de Riva%NULL%1,                Ana%Ruiz%NULL%1,                Beatriz%Tena%NULL%1,                Sebastián%Jaramillo%NULL%1,                José María%Balibrea%NULL%1,                Francisco Borja de Borja%de Lacy%NULL%1,                Ana%Otero%NULL%1,                Ainitze%Ibarzabal%NULL%1,                Raquel%Bravo%NULL%1,                Anna%Carreras%NULL%1,                Daniel%Martín-Barreda%NULL%1,                Alfonso Jesús%Alias%NULL%1,                Mariano%Balaguer%NULL%1,                Jorge%Aliaga%NULL%1,                Alex%Almuedo%NULL%1,                Joan Ramón%Alonso%NULL%1,                Rut%Andrea%NULL%1,                Gerard Sergi%Angelès%NULL%1,                 Marilyn%Arias%NULL%1,                Fátima%Aziz%NULL%1,                Joan Ramon%Badía%NULL%1,                Enric%Barbeta%NULL%1,                Toni%Torres%NULL%1,                Guillem%Batiste%NULL%1,                Pau%Benet%NULL%1,                Xavi%Borrat%NULL%1,                María%Borrell%NULL%1,                Ernest%Bragulat%NULL%1,                Inmaculada%Carmona%NULL%1,                Manuel%Castellà%NULL%1,                Pedro%Castro%NULL%1,                Joan%Ceravalls%NULL%1,                Oscar%Comino%NULL%1,                Claudia%Cucciniello%NULL%1,                Clàudia%De Deray%NULL%1,                Oriol%De Diego%NULL%1,                Paula%De la Matta%NULL%1,                Marta%Farrero%NULL%1,                Javier%Fernández%NULL%1,                Sara%Fernández%NULL%1,                Anna%Fernández%NULL%1,                Miquel%Ferrer%NULL%1,                Ana%Fervienza%NULL%1,                María%Tallo Forga%NULL%1,                Daniel%Forné%NULL%1,                Clàudia%Galán%NULL%1,                Andrea%Gómez%NULL%1,                Eduard%Guasch%NULL%1,                María%Hernández-Tejero%NULL%1,                Adriana%Jacas%NULL%1,                Beltrán%Jiménez%NULL%1,                Pere%Leyes%NULL%1,                Teresa%López%NULL%1,                José Antonio%Martínez%NULL%1,                Graciela%Martínez-Pallí%NULL%1,                Jordi%Mercadal%NULL%1,                Guido%Muñoz%NULL%1,                José%Muñoz%NULL%1,                Ricard%Navarro%NULL%0,                Josep María%Nicolás%NULL%1,                José Tomás%Ortiz%NULL%1,                Anna%Peiró%NULL%1,                Manuel%Pérez%NULL%1,                Esteban%Poch%NULL%1,                Margarida%Pujol%NULL%1,                Eduard%Quintana%NULL%1,                Bartomeu%Ramis%NULL%1,                Enric%Reverter%NULL%1,                Irene%Rovira%NULL%1,                Pablo%Ruiz%NULL%1,                Elena%Sandoval%NULL%1,                Stefan%Schneider%NULL%1,                Oriol%Sibila%NULL%1,                Carla%Solé%NULL%1,                Alex%Soriano%NULL%1,                Dolors%Soy%NULL%1,                M.%Suárez%NULL%1,                Adrián%Téllez%NULL%1,                Néstor David%Toapanta%NULL%1,                Antoni%Torres%NULL%0,                Xavier%Urra%NULL%1,                César%Aldecoa%NULL%1,                Alicia%Bordell%NULL%1,                Silvia%Martín%NULL%1,                Judith%Andrés%NULL%1,                Alberto Martínez%Ruiz%NULL%1,                Gonzalo Tamayo%Medel%NULL%1,                Iñaki Bilbao Bilbao%Villasante%NULL%1,                Fernando Iturri%Clavero%NULL%1,                Covadonga Peralta%Álvarez%NULL%1,                Julia T.%Herrera Díez%NULL%1,                Andrea García%Trancho%NULL%1,                Iñaki Sainz %Mandiola%NULL%1,                Carmen Ruano%Suarez%NULL%1,                Angela Ruiz%Bocos%NULL%1,                Eneritz Urrutia%Izagirre%NULL%1,                Pablo Ortiz%de Urbina Fernández%NULL%1,                Naiara Apodaka%López%NULL%1,                Leire Prieto%Molano%NULL%1,                Eunate Ganuza%Martínez%NULL%1,                Iratxe Vallinas%Hidalgo%NULL%1,                Karmele%de Orte Sancho%NULL%1,                Celia González%Paniagua%NULL%1,                Gemma Ortiz%Labrador%NULL%1,                Mireia Pérez%Larrañaga%NULL%1,                Marta López%Miguelez%NULL%1,                Estíbaliz Bárcena%Andrés%NULL%1,                Erik Urutxurtu%Laureano%NULL%1,                Maria Jesús Maroño%Boedo%NULL%1,                Blanca Escontrela%Rodríguez%NULL%1,                Aitziber Ereñozaga%Camiruaga%NULL%1,                Deiene Lasuen%Aguirre%NULL%1,                Ainhoa Zabal%Maeztu%NULL%1,                Ane Guereca%Gala%NULL%1,                Iker Castelo%Korro%NULL%1,                Andrés Álvarez%Campo%NULL%1,                Alejandro Carcelen%Viana%NULL%1,                Alejandro Alberdi%Enríquez%NULL%1,                Xabier Ormazábal%Rementeria%NULL%1,                Alberto Sánchez%Campos%NULL%1,                Rosa Gutiérrez%Rico%NULL%1,                Pablo Barbier%Damborenea%NULL%1,                Marta Guerenabarrena%Momeñe%NULL%1,                Borja Cuesta%Ruiz%NULL%1,                Alejandro López%Rico%NULL%1,                Ana Rojo%Polo%NULL%1,                Covadonga García%Grijelmo%NULL%1,                Mikel Celorrio%Reta%NULL%1,                Eneko Martín%Arroyo%NULL%1,                Leire Artaza%Aparicio%NULL%1,                Iñaki Ituarte%Aspiazu%NULL%1,                Ane Igeregi%Basabe%NULL%1,                Itxaso Merino%Julian%NULL%1,                Isabel Diaz%Rico%NULL%1,                Maria Paz%Martínez%NULL%1,                Ramón%Adalia Bartolomé%NULL%1,                Luigi%Zattera%NULL%1,                Irina Adalid%Hernandez%NULL%1,                Leire Larrañaga%Altuna%NULL%1,                Aina Serrallonga%Castells%NULL%1,                Adriana Vílchez%Garcia%NULL%1,                María%Núñez%NULL%1,                Lorena%Román%NULL%1,                Isabel Ramos%Delgado%NULL%1,                Adela Benítez-Cano%Martínez%NULL%1,                Mireia Chanzá%Albert%NULL%1,                Juan Carlos Álvarez%García%NULL%1,                Luis Aguilera%Cuchillo%NULL%1,                Sandra Beltrán%de Heredia%NULL%1,                Jesús Carazo%Cordobés%NULL%1,                Carlos Alberto García%Bernedo%NULL%1,                Fernando Escolano%Villén%NULL%1,                Francisco Javier Redondo%Calvo%NULL%1,                Rubén Villazala%González%NULL%1,                Victor Baladron%González%NULL%1,                Patricia%Faba%NULL%1,                Omar%Montenegro%NULL%1,                Natalia Bejarano%Ramírez%NULL%1,                Sergio Marcos%Contreras%NULL%1,                Alejandro Garcia%Rodríguez%NULL%1,                Saleta Rey%Vázquez%NULL%1,                Cristina Garcia%Pérez%NULL%1,                Eva Higuera%Miguelez%NULL%1,                Irene Pérez%Blanco%NULL%1,                David García%Rivera%NULL%1,                Ane Martín%de la Fuente%NULL%1,                Marta%Pardo%NULL%1,                Vanessa%Rodriguez%NULL%1,                Unai%Bengoetxea%NULL%1,                Fernando%Ramasco%NULL%1,                Sheila Olga Santidrián%Bernal%NULL%1,                Alvar Santa Cruz%Hernando%NULL%1,                Antonio Planas%Roca%NULL%1,                Carlos Figueroa%Yusta%NULL%1,                Esther García%Villabona%NULL%1,                Carmen Vallejo%Lantero%NULL%1,                Eva Patiño%Rodriguez%NULL%1,                Alvaro Esquivel%Toledo%NULL%1,                David Arribas%Méndez%NULL%1,                Mar Orts%Rodriguez%NULL%1,                Rosa Méndez%Hernández%NULL%1,                Jesús Nieves%Alonso%NULL%1,                Inés Imaz%Artazcoz%NULL%1,                Sonia Expósito%Carazo%NULL%1,                Carlos Román%Guerrero%NULL%1,                Elena Rojo%Rodríguez%NULL%1,                Ricardo Moreno%González%NULL%1,                Julia Hernando%Santos%NULL%1,                Jara Torrente%Pérez%NULL%1,                Esperanza Mata%Mena%NULL%1,                Manuel José Muñoz%Martínez%NULL%1,                Enrique Alday%Muñoz%NULL%1,                Patricia Martin%Serrano%NULL%1,                Laura Cotter%Muñoz%NULL%1,                Amadea%Mjertan%NULL%1,                Diego Gutierrez%Martínez%NULL%1,                Carmen Rodríguez%García%NULL%1,                Olaya Alonso%Viejo%NULL%1,                Juan Alvarez%Pereira%NULL%1,                Ana Carmona%Bonet%NULL%1,                Diana Parrado%López%NULL%1,                Eva%de Dios Tomas%NULL%1,                Rafael Martín%Celemin%NULL%1,                María Luisa Meilan%Paz%NULL%1,                Luis Quecedo%Gutiérrez%NULL%1,                Noemí Diaz%Velasco%NULL%1,                Gabriel Martin%Hernández%NULL%1,                Francisco Garcia%del Corral%NULL%1,                Gloria Hernandez%Arias%NULL%1,                David Rodriguez%Cuesta%NULL%1,                Ana Gómez%Rice%NULL%1,                Encarna Mateos%Sevillano%NULL%1,                Natalia Olmos%Molpeceres%NULL%1,                Beatriz%Domínguez%NULL%1,                Ana Vázquez%Lima%NULL%1,                Ángel%Candela%NULL%1,                Ismael A. Acevedo%Bambaren%NULL%1,                Maria Isabel Albala%Blanco%NULL%1,                Paloma Alonso%Montoiro%NULL%1,                Fernando Álvarez%Utrera%NULL%1,                Juan Avellanosa%Esteruelas%NULL%1,                Amal Azzam%López%NULL%1,                Alberto José%Balvis%NULL%1,                Balvis Tommaso%Bardi%NULL%1,                María Beltrán%Martín%NULL%1,                Jacobo Benatar%Haserfaty%NULL%1,                Alberto Berruezo%Camacho%NULL%1,                Laura Betolaza%Weimer%NULL%1,                María%del Mar Carbonell Soto%NULL%1,                Cristina Carrasco%Seral%NULL%1,                Cristina Cerro %Zaballos%NULL%1,                Elizabeth Claros%Llamas%NULL%1,                Pilar Coleta%Orduna%NULL%1,                Ingrid P. Cortes%Forero%NULL%1,                Pascual Agustín Crespo%Aliseda%NULL%1,                María Angélica%de Pablo Pajares%NULL%1,                Yolanda Díez%Remesal%NULL%1,                Trinidad Dorado%Díaz%NULL%1,                Noemí Echevarría%Blasco%NULL%1,                María Elena Elías%Martín%NULL%1,                Javier Felices%Triviño%NULL%1,                Natalia Fernández%López%NULL%1,                Cristina Fernández%Martín%NULL%1,                Natalia Ferreiro%Pozuelo%NULL%1,                Luis Gajate%Martín%NULL%1,                Clara Gallego%Santos%NULL%1,                Diego Gil%Mayo%NULL%1,                María Gómez%Rojo%NULL%1,                Claudia González%Cibrián%NULL%1,                Elena Herrera%López%NULL%1,                Borja Hinojal%Olmedillo%NULL%1,                Berta Iglesias%Gallego%NULL%1,                Sassan%Khonsari%NULL%1,                María Nuria Mane%Ruiz%NULL%1,                María Manzanero%Arroyo%NULL%1,                Ana María Mariscal%Ortega%NULL%1,                Sara Martín%Burcio%NULL%1,                María%del Carmen Martín González%NULL%1,                Ascensión Martín%Grande%NULL%1,                Jose Juan Martín%López%NULL%1,                Cecilia Martín%Rabes%NULL%1,                Marcos Martínez%Borja%NULL%1,                Nilda Martínez%Castro%NULL%1,                Adolfo Martínez%Pérez%NULL%1,                Snejana%Matcan%NULL%1,                Cristina Medrano%Viñas%NULL%1,                Lisset Miguel%Herrera%NULL%1,                Adrián Mira%Betancur%NULL%1,                María Montiel%Carbajo%NULL%1,                Javier Moya%Moradas%NULL%1,                Lorena Muñoz%Pérez%NULL%1,                Mónica Nuñez%Murias%NULL%1,                Eva Ordiales%González%NULL%1,                Óscar Ordoñez%Recio%NULL%1,                Miguel Ángel Palomero%Rodriguez%NULL%1,                Diego Parise%Roux%NULL%1,                Lucia Pereira%Torres%NULL%1,                David Pestaña%Lagunas%NULL%1,                Juana María Pinto%Corraliza%NULL%1,                Marian Prieto%Rodrigo%NULL%1,                Inmaculada Rodriguez%Diaz-Regaño%NULL%1,                David Rodriguez%Esteban%NULL%1,                Víctor Rojas%Pernia%NULL%1,                Álvaro Ruigómez%Saiz%NULL%1,                Bárbara Saavedra%Villarino%NULL%1,                Noemí Samaranch%Palero%NULL%1,                Gloria Santos%Pérez%NULL%1,                Jaume Serna%Pérez%NULL%1,                Ana Belén Serrano%Romero%NULL%1,                Jesús Tercero%López%NULL%1,                Carlos Tiscar%García%NULL%1,                Marta%de la Torre Concostrina%NULL%1,                Eva María Ureta%Mesa%NULL%1,                Eva Velasco%Olarte%NULL%1,                Judith Villahoz%Martínez%NULL%1,                Raúl Villalaba%Palacios%NULL%1,                Gema Villanueva%García%NULL%1,                Cristina Vogel%de Medeiros%NULL%1,                Soraya Gholamian%Ovejero%NULL%1,                Marta Vicente%Orgaz%NULL%1,                Patricia Lloreda%Herradon%NULL%1,                Cristina Crespo%Gómez%NULL%1,                Tatiana%Sarmiento-Trujillo%NULL%1,                Noemí García%Medina%NULL%1,                María Martínez%García%NULL%1,                Carles Espinós%Ramírez%NULL%1,                Nabil Mouhaffel%Rivero%NULL%1,                Jose Antonio Bernia%Gil%NULL%1,                Sonsoles%Martín%NULL%1,                María Victoria%Moral%NULL%1,                Josefina%Galán%NULL%1,                Pilar%Paniagua%NULL%1,                Sergio%Pérez%NULL%1,                Albert%Bainac%NULL%1,                Ana%Arias%NULL%1,                Elsa%Ramil%NULL%1,                Jorge%Escudero%NULL%1,                Pablo%Monedero%NULL%1,                Carmen%Cara%NULL%1,                Andrea%Lara%NULL%1,                Elena Mendez%Martínez%NULL%1,                Jorge%Mendoza%NULL%1,                Íñigo Rubio%Baines%NULL%1,                Carmen Sala%Trull%NULL%1,                Pablo Montero%López%NULL%1,                Alfredo%Gea%NULL%1,                Alejandro%Montero%NULL%1,                Rocío Armero%Ibañez%NULL%1,                Juan Vicente Llau%Pitarch%NULL%1,                Fernando Rauer%Alcóver%NULL%1,                Cristina Álvarez%Herreros%NULL%1,                Cyntia Sánchez%Martín%NULL%1,                Lucía López Ocáriz%Olmos%NULL%1,                Marta Navas%Moruno%NULL%1,                Fernando García%Montoto%NULL%1,                M. F. Mirón%Rodriguez%NULL%1,                Laura Fuentes%Coco%NULL%1,                Cristina Hernández%Gamito%NULL%1,                Antonio Barba%Orejudo%NULL%1,                Luis Gerardo Smith%Vielma%NULL%1,                Yasmina González Marín Francisco%de Borja Amador Penco%NULL%1,                Marta Donoso%Domínguez%NULL%1,                Silvia Esquivel%Ramírez%NULL%1,                José Antonio%Carbonell%NULL%1,                Berta Monleón%López%NULL%1,                Sara%Martínez-Castro%NULL%1,                Gerardo%Aguilar%NULL%1,                María Gestal Pablo%Casas%NULL%1,                Angel Outeiro%Rosato%NULL%1,                Andrea Naveiro%Pan%NULL%1,                María Alonso%Portela%NULL%1,                Adrián García%Romar%NULL%1,                Eva Mosquera%Rodríguez%NULL%1,                Diego Ruanova%Seijo%NULL%1,                Pablo Rama%Maceiras%NULL%1,                Francisco%Castro-Ceoane%NULL%1,                Esther Moreno%López%NULL%1,                Sergio%Gil%NULL%1,                Julia Guillén%Antón%NULL%1,                Patricia García-Consuegra%Tirado%NULL%1,                Aurora Callau%Calvo%NULL%1,                Laura Forés%Lisbona%NULL%1,                María Carbonell%Romero%NULL%1,                Belén Albericio%Gil%NULL%1,                Laura Pradal%Jarne%NULL%1,                María Soria%Lozano%NULL%1,                Diego Loscos%López%NULL%1]'
$ws.Range("E12").Value2 = '[Aloknath%Pandya%NULL%1,                Navjot Ariyana%Kaur%NULL%1,                Daniel%Sacher%NULL%1,                Oisin%O’Corragain%NULL%1,                Daniel%Salerno%NULL%1,                Parag%Desai%NULL%1,                Sameep%Sehgal%NULL%1,                Matthew%Gordon%NULL%2,                Rohit%Gupta%NULL%2,                Nathaniel%Marchetti%NULL%1,                Huaqing%Zhao%NULL%2,                Nicole%Patlakh%NULL%2,                Gerard J.%Criner%NULL%1,                Temple%University%NULL%1,                NULL%NULL%NULL%0]'
$ws.Range("E13").Value2 = '[Kevin%Roedl%NULL%2,                Dominik%Jarczak%NULL%1,                Liina%Thasler%NULL%1,                Martin%Bachmann%NULL%1,                Frank%Schulte%NULL%1,                Berthold%Bein%NULL%1,                Christian Friedrich%Weber%NULL%1,                Ulrich%Schäfer%NULL%1,                Carsten%Veit%NULL%1,                Hans-Peter%Hauber%NULL%1,                Sebastian%Kopp%NULL%1,                Karsten%Sydow%NULL%1,                Andreas%de Weerth%NULL%0,                Marc%Bota%NULL%1,                Rüdiger%Schreiber%NULL%1,                Oliver%Detsch%NULL%1,                Jan-Peer%Rogmann%NULL%1,                Daniel%Frings%NULL%0,                Barbara%Sensen%NULL%1,                Christoph%Burdelski%NULL%0,                Olaf%Boenisch%NULL%1,                Axel%Nierhaus%NULL%0,                Geraldine%de Heer%NULL%0,                Stefan%Kluge%NULL%0]'
